$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 42.1867505  # G2
$ws.Cells.Item(2, 8).Value = 84.373501  # H2
$ws.Cells.Item(2, 9).Value = 0.36721878565825766  # I2
$ws.Cells.Item(2, 10).Value = 0.30197077016319174  # J2
$ws.Cells.Item(2, 13).Value = 228.691925  # M2
$ws.Cells.Item(2, 14).Value = 457.38385  # N2
$ws.Cells.Item(2, 15).Value = 0.5727940074410814  # O2
$ws.Cells.Item(2, 16).Value = 0.475934002654817  # P2
$ws.Cells.Item(2, 17).Value = 9647.769181339712  # Q2
$ws.Cells.Item(2, 18).Value = 38591.07672535885  # R2
$ws.Cells.Item(2, 19).Value = 0.21034071984484093  # S2
$ws.Cells.Item(2, 20).Value = 0.14371815732852564  # T2
$ws.Cells.Item(3, 7).Value = 42.1867505  # G3
$ws.Cells.Item(3, 8).Value = 84.373501  # H3
$ws.Cells.Item(3, 9).Value = 0.36721878565825766  # I3
$ws.Cells.Item(3, 10).Value = 0.30197077016319174  # J3
$ws.Cells.Item(3, 15).Value = 0.0011895817926618896  # O3
$ws.Cells.Item(3, 16).Value = 0.0014826335909033635  # P3
$ws.Cells.Item(3, 17).Value = 20.036540900974334  # Q3
$ws.Cells.Item(3, 18).Value = 120.21924540584601  # R3
$ws.Cells.Item(3, 19).Value = 0.00043683678134247235  # S3
$ws.Cells.Item(3, 20).Value = 0.0004477120073149072  # T3
$ws.Cells.Item(4, 7).Value = 42.1867505  # G4
$ws.Cells.Item(4, 8).Value = 84.373501  # H4
$ws.Cells.Item(4, 9).Value = 0.36721878565825766  # I4
$ws.Cells.Item(4, 10).Value = 0.30197077016319174  # J4
$ws.Cells.Item(4, 13).Value = 38.30460633333333  # M4
$ws.Cells.Item(4, 14).Value = 114.913819  # N4
$ws.Cells.Item(4, 15).Value = 0.0959397624779407  # O4
$ws.Cells.Item(4, 16).Value = 0.11957438776428411  # P4
$ws.Cells.Item(4, 17).Value = 1615.9468703850532  # Q4
$ws.Cells.Item(4, 18).Value = 9695.68122231032  # R4
$ws.Cells.Item(4, 19).Value = 0.03523088307349106  # S4
$ws.Cells.Item(4, 20).Value = 0.036107969964973  # T4
$ws.Cells.Item(5, 7).Value = 42.1867505  # G5
$ws.Cells.Item(5, 8).Value = 84.373501  # H5
$ws.Cells.Item(5, 9).Value = 0.36721878565825766  # I5
$ws.Cells.Item(5, 10).Value = 0.30197077016319174  # J5
$ws.Cells.Item(5, 13).Value = 8.054906  # M5
$ws.Cells.Item(5, 14).Value = 16.109812  # N5
$ws.Cells.Item(5, 15).Value = 0.02017474769737153  # O5
$ws.Cells.Item(5, 16).Value = 0.01676317890799293  # P5
$ws.Cells.Item(5, 17).Value = 339.81030972295304  # Q5
$ws.Cells.Item(5, 18).Value = 1359.2412388918121  # R5
$ws.Cells.Item(5, 19).Value = 0.007408546350390502  # S5
$ws.Cells.Item(5, 20).Value = 0.0050619900452299966  # T5
$ws.Cells.Item(6, 7).Value = 42.1867505  # G6
$ws.Cells.Item(6, 8).Value = 84.373501  # H6
$ws.Cells.Item(6, 9).Value = 0.36721878565825766  # I6
$ws.Cells.Item(6, 10).Value = 0.30197077016319174  # J6
$ws.Cells.Item(6, 13).Value = 21.28819633333333  # M6
$ws.Cells.Item(6, 14).Value = 63.864589  # N6
$ws.Cells.Item(6, 15).Value = 0.05331955332031306  # O6
$ws.Cells.Item(6, 16).Value = 0.06645475014186618  # P6
$ws.Cells.Item(6, 17).Value = 898.0798273093482  # Q6
$ws.Cells.Item(6, 18).Value = 5388.47896385609  # R6
$ws.Cells.Item(6, 19).Value = 0.019579941622126083  # S6
$ws.Cells.Item(6, 20).Value = 0.020067392081341805  # T6
$ws.Cells.Item(7, 7).Value = 42.1867505  # G7
$ws.Cells.Item(7, 8).Value = 84.373501  # H7
$ws.Cells.Item(7, 9).Value = 0.36721878565825766  # I7
$ws.Cells.Item(7, 10).Value = 0.30197077016319174  # J7
$ws.Cells.Item(7, 13).Value = 102.442257  # M7
$ws.Cells.Item(7, 14).Value = 307.326771  # N7
$ws.Cells.Item(7, 15).Value = 0.25658234727063134  # O7
$ws.Cells.Item(7, 16).Value = 0.31979104694013644  # P7
$ws.Cells.Item(7, 17).Value = 4321.705936715879  # Q7
$ws.Cells.Item(7, 18).Value = 25930.235620295272  # R7
$ws.Cells.Item(7, 19).Value = 0.09422185798606661  # S7
$ws.Cells.Item(7, 20).Value = 0.0965675487358064  # T7
$ws.Cells.Item(8, 9).Value = 0.4297349423931389  # I8
$ws.Cells.Item(8, 10).Value = 0.5300684357196348  # J8
$ws.Cells.Item(8, 13).Value = 228.691925  # M8
$ws.Cells.Item(8, 14).Value = 457.38385  # N8
$ws.Cells.Item(8, 15).Value = 0.5727940074410814  # O8
$ws.Cells.Item(8, 16).Value = 0.475934002654817  # P8
$ws.Cells.Item(8, 17).Value = 11290.22723043279  # Q8
$ws.Cells.Item(8, 18).Value = 67741.36338259676  # R8
$ws.Cells.Item(8, 19).Value = 0.24614959979082832  # S8
$ws.Cells.Item(8, 20).Value = 0.2522775922930234  # T8
$ws.Cells.Item(9, 9).Value = 0.4297349423931389  # I9
$ws.Cells.Item(9, 10).Value = 0.5300684357196348  # J9
$ws.Cells.Item(9, 15).Value = 0.0011895817926618896  # O9
$ws.Cells.Item(9, 16).Value = 0.0014826335909033635  # P9
$ws.Cells.Item(9, 19).Value = 0.0005112048631414841  # S9
$ws.Cells.Item(9, 20).Value = 0.0007858972682755308  # T9
$ws.Cells.Item(10, 9).Value = 0.4297349423931389  # I10
$ws.Cells.Item(10, 10).Value = 0.5300684357196348  # J10
$ws.Cells.Item(10, 13).Value = 38.30460633333333  # M10
$ws.Cells.Item(10, 14).Value = 114.913819  # N10
$ws.Cells.Item(10, 15).Value = 0.0959397624779407  # O10
$ws.Cells.Item(10, 16).Value = 0.11957438776428411  # P10
$ws.Cells.Item(10, 17).Value = 1891.049320939549  # Q10
$ws.Cells.Item(10, 18).Value = 17019.443888455946  # R10
$ws.Cells.Item(10, 19).Value = 0.041228668301669276  # S10
$ws.Cells.Item(10, 20).Value = 0.06338260867434711  # T10
$ws.Cells.Item(11, 9).Value = 0.4297349423931389  # I11
$ws.Cells.Item(11, 10).Value = 0.5300684357196348  # J11
$ws.Cells.Item(11, 13).Value = 8.054906  # M11
$ws.Cells.Item(11, 14).Value = 16.109812  # N11
$ws.Cells.Item(11, 15).Value = 0.02017474769737153  # O11
$ws.Cells.Item(11, 16).Value = 0.01676317890799293  # P11
$ws.Cells.Item(11, 17).Value = 397.66038551547666  # Q11
$ws.Cells.Item(11, 18).Value = 2385.9623130928603  # R11
$ws.Cells.Item(11, 19).Value = 0.008669794039526066  # S11
$ws.Cells.Item(11, 20).Value = 0.00888563202144819  # T11
$ws.Cells.Item(12, 9).Value = 0.4297349423931389  # I12
$ws.Cells.Item(12, 10).Value = 0.5300684357196348  # J12
$ws.Cells.Item(12, 13).Value = 21.28819633333333  # M12
$ws.Cells.Item(12, 14).Value = 63.864589  # N12
$ws.Cells.Item(12, 15).Value = 0.05331955332031306  # O12
$ws.Cells.Item(12, 16).Value = 0.06645475014186618  # P12
$ws.Cells.Item(12, 17).Value = 1050.9709686050326  # Q12
$ws.Cells.Item(12, 18).Value = 9458.738717445296  # R12
$ws.Cells.Item(12, 19).Value = 0.02291327517453263  # S12
$ws.Cells.Item(12, 20).Value = 0.03522556545383818  # T12
$ws.Cells.Item(13, 9).Value = 0.4297349423931389  # I13
$ws.Cells.Item(13, 10).Value = 0.5300684357196348  # J13
$ws.Cells.Item(13, 13).Value = 102.442257  # M13
$ws.Cells.Item(13, 14).Value = 307.326771  # N13
$ws.Cells.Item(13, 15).Value = 0.25658234727063134  # O13
$ws.Cells.Item(13, 16).Value = 0.31979104694013644  # P13
$ws.Cells.Item(13, 17).Value = 5057.442931263944  # Q13
$ws.Cells.Item(13, 18).Value = 45516.98638137551  # R13
$ws.Cells.Item(13, 19).Value = 0.11026240022344112  # S13
$ws.Cells.Item(13, 20).Value = 0.16951114000870243  # T13
$ws.Cells.Item(14, 5).Value = 1.0  # E14
$ws.Cells.Item(14, 6).Value = 0.3333333333333333  # F14
$ws.Cells.Item(14, 7).Value = 0.006662  # G14
$ws.Cells.Item(14, 8).Value = 0.019986  # H14
$ws.Cells.Item(14, 9).Value = 0.000057990044766669394  # I14
$ws.Cells.Item(14, 10).Value = 0.0000715294226380573  # J14
$ws.Cells.Item(14, 13).Value = 228.691925  # M14
$ws.Cells.Item(14, 14).Value = 457.38385  # N14
$ws.Cells.Item(14, 15).Value = 0.5727940074410814  # O14
$ws.Cells.Item(14, 16).Value = 0.475934002654817  # P14
$ws.Cells.Item(14, 17).Value = 1.52354560435  # Q14
$ws.Cells.Item(14, 18).Value = 9.1412736261  # R14
$ws.Cells.Item(14, 19).Value = 0.000033216350133588276  # S14
$ws.Cells.Item(14, 20).Value = 0.00003404328442371869  # T14
$ws.Cells.Item(15, 5).Value = 1.0  # E15
$ws.Cells.Item(15, 6).Value = 0.3333333333333333  # F15
$ws.Cells.Item(15, 7).Value = 0.006662  # G15
$ws.Cells.Item(15, 8).Value = 0.019986  # H15
$ws.Cells.Item(15, 9).Value = 0.000057990044766669394  # I15
$ws.Cells.Item(15, 10).Value = 0.0000715294226380573  # J15
$ws.Cells.Item(15, 15).Value = 0.0011895817926618896  # O15
$ws.Cells.Item(15, 16).Value = 0.0014826335909033635  # P15
$ws.Cells.Item(15, 17).Value = 0.0031641080173333337  # Q15
$ws.Cells.Item(15, 18).Value = 0.028476972156  # R15
$ws.Cells.Item(15, 19).Value = 0.00000006898390141007781  # S15
$ws.Cells.Item(15, 20).Value = 0.00000010605192474110722  # T15
$ws.Cells.Item(16, 5).Value = 1.0  # E16
$ws.Cells.Item(16, 6).Value = 0.3333333333333333  # F16
$ws.Cells.Item(16, 7).Value = 0.006662  # G16
$ws.Cells.Item(16, 8).Value = 0.019986  # H16
$ws.Cells.Item(16, 9).Value = 0.000057990044766669394  # I16
$ws.Cells.Item(16, 10).Value = 0.0000715294226380573  # J16
$ws.Cells.Item(16, 13).Value = 38.30460633333333  # M16
$ws.Cells.Item(16, 14).Value = 114.913819  # N16
$ws.Cells.Item(16, 15).Value = 0.0959397624779407  # O16
$ws.Cells.Item(16, 16).Value = 0.11957438776428411  # P16
$ws.Cells.Item(16, 17).Value = 0.25518528739266666  # Q16
$ws.Cells.Item(16, 18).Value = 2.2966675865340003  # R16
$ws.Cells.Item(16, 19).Value = 0.00000556355112099941  # S16
$ws.Cells.Item(16, 20).Value = 0.000008553086919078425  # T16
$ws.Cells.Item(17, 5).Value = 1.0  # E17
$ws.Cells.Item(17, 6).Value = 0.3333333333333333  # F17
$ws.Cells.Item(17, 7).Value = 0.006662  # G17
$ws.Cells.Item(17, 8).Value = 0.019986  # H17
$ws.Cells.Item(17, 9).Value = 0.000057990044766669394  # I17
$ws.Cells.Item(17, 10).Value = 0.0000715294226380573  # J17
$ws.Cells.Item(17, 13).Value = 8.054906  # M17
$ws.Cells.Item(17, 14).Value = 16.109812  # N17
$ws.Cells.Item(17, 15).Value = 0.02017474769737153  # O17
$ws.Cells.Item(17, 16).Value = 0.01676317890799293  # P17
$ws.Cells.Item(17, 17).Value = 0.05366178377200001  # Q17
$ws.Cells.Item(17, 18).Value = 0.32197070263200006  # R17
$ws.Cells.Item(17, 19).Value = 0.0000011699345221268352  # S17
$ws.Cells.Item(17, 20).Value = 0.0000011990605088671942  # T17
$ws.Cells.Item(18, 5).Value = 1.0  # E18
$ws.Cells.Item(18, 6).Value = 0.3333333333333333  # F18
$ws.Cells.Item(18, 7).Value = 0.006662  # G18
$ws.Cells.Item(18, 8).Value = 0.019986  # H18
$ws.Cells.Item(18, 9).Value = 0.000057990044766669394  # I18
$ws.Cells.Item(18, 10).Value = 0.0000715294226380573  # J18
$ws.Cells.Item(18, 13).Value = 21.28819633333333  # M18
$ws.Cells.Item(18, 14).Value = 63.864589  # N18
$ws.Cells.Item(18, 15).Value = 0.05331955332031306  # O18
$ws.Cells.Item(18, 16).Value = 0.06645475014186618  # P18
$ws.Cells.Item(18, 17).Value = 0.14182196397266666  # Q18
$ws.Cells.Item(18, 18).Value = 1.276397675754  # R18
$ws.Cells.Item(18, 19).Value = 0.0000030920032839837703  # S18
$ws.Cells.Item(18, 20).Value = 0.000004753469909204044  # T18
$ws.Cells.Item(19, 5).Value = 1.0  # E19
$ws.Cells.Item(19, 6).Value = 0.3333333333333333  # F19
$ws.Cells.Item(19, 7).Value = 0.006662  # G19
$ws.Cells.Item(19, 8).Value = 0.019986  # H19
$ws.Cells.Item(19, 9).Value = 0.000057990044766669394  # I19
$ws.Cells.Item(19, 10).Value = 0.0000715294226380573  # J19
$ws.Cells.Item(19, 13).Value = 102.442257  # M19
$ws.Cells.Item(19, 14).Value = 307.326771  # N19
$ws.Cells.Item(19, 15).Value = 0.25658234727063134  # O19
$ws.Cells.Item(19, 16).Value = 0.31979104694013644  # P19
$ws.Cells.Item(19, 17).Value = 0.682470316134  # Q19
$ws.Cells.Item(19, 18).Value = 6.142232845206  # R19
$ws.Cells.Item(19, 19).Value = 0.000014879221804561023  # S19
$ws.Cells.Item(19, 20).Value = 0.00002287446895244784  # T19
$ws.Cells.Item(20, 7).Value = 23.049115  # G20
$ws.Cells.Item(20, 8).Value = 46.09823  # H20
$ws.Cells.Item(20, 9).Value = 0.2006333249297675  # I20
$ws.Cells.Item(20, 10).Value = 0.16498447796139157  # J20
$ws.Cells.Item(20, 13).Value = 228.691925  # M20
$ws.Cells.Item(20, 14).Value = 457.38385  # N20
$ws.Cells.Item(20, 15).Value = 0.5727940074410814  # O20
$ws.Cells.Item(20, 16).Value = 0.475934002654817  # P20
$ws.Cells.Item(20, 17).Value = 5271.146478896375  # Q20
$ws.Cells.Item(20, 18).Value = 21084.5859155855  # R20
$ws.Cells.Item(20, 19).Value = 0.11492156621275015  # S20
$ws.Cells.Item(20, 20).Value = 0.07852172297208053  # T20
$ws.Cells.Item(21, 7).Value = 23.049115  # G21
$ws.Cells.Item(21, 8).Value = 46.09823  # H21
$ws.Cells.Item(21, 9).Value = 0.2006333249297675  # I21
$ws.Cells.Item(21, 10).Value = 0.16498447796139157  # J21
$ws.Cells.Item(21, 15).Value = 0.0011895817926618896  # O21
$ws.Cells.Item(21, 16).Value = 0.0014826335909033635  # P21
$ws.Cells.Item(21, 17).Value = 10.947146437096666  # Q21
$ws.Cells.Item(21, 18).Value = 65.68287862258  # R21
$ws.Cells.Item(21, 19).Value = 0.0002386697503376682  # S21
$ws.Cells.Item(21, 20).Value = 0.0002446115290032148  # T21
$ws.Cells.Item(22, 7).Value = 23.049115  # G22
$ws.Cells.Item(22, 8).Value = 46.09823  # H22
$ws.Cells.Item(22, 9).Value = 0.2006333249297675  # I22
$ws.Cells.Item(22, 10).Value = 0.16498447796139157  # J22
$ws.Cells.Item(22, 13).Value = 38.30460633333333  # M22
$ws.Cells.Item(22, 14).Value = 114.913819  # N22
$ws.Cells.Item(22, 15).Value = 0.0959397624779407  # O22
$ws.Cells.Item(22, 16).Value = 0.11957438776428411  # P22
$ws.Cells.Item(22, 17).Value = 882.8872764067283  # Q22
$ws.Cells.Item(22, 18).Value = 5297.32365844037  # R22
$ws.Cells.Item(22, 19).Value = 0.01924871353892139  # S22
$ws.Cells.Item(22, 20).Value = 0.019727917942843422  # T22
$ws.Cells.Item(23, 7).Value = 23.049115  # G23
$ws.Cells.Item(23, 8).Value = 46.09823  # H23
$ws.Cells.Item(23, 9).Value = 0.2006333249297675  # I23
$ws.Cells.Item(23, 10).Value = 0.16498447796139157  # J23
$ws.Cells.Item(23, 13).Value = 8.054906  # M23
$ws.Cells.Item(23, 14).Value = 16.109812  # N23
$ws.Cells.Item(23, 15).Value = 0.02017474769737153  # O23
$ws.Cells.Item(23, 16).Value = 0.01676317890799293  # P23
$ws.Cells.Item(23, 17).Value = 185.65845470819002  # Q23
$ws.Cells.Item(23, 18).Value = 742.6338188327601  # R23
$ws.Cells.Item(23, 19).Value = 0.00404772671014282  # S23
$ws.Cells.Item(23, 20).Value = 0.0027656643211086238  # T23
$ws.Cells.Item(24, 7).Value = 23.049115  # G24
$ws.Cells.Item(24, 8).Value = 46.09823  # H24
$ws.Cells.Item(24, 9).Value = 0.2006333249297675  # I24
$ws.Cells.Item(24, 10).Value = 0.16498447796139157  # J24
$ws.Cells.Item(24, 13).Value = 21.28819633333333  # M24
$ws.Cells.Item(24, 14).Value = 63.864589  # N24
$ws.Cells.Item(24, 15).Value = 0.05331955332031306  # O24
$ws.Cells.Item(24, 16).Value = 0.06645475014186618  # P24
$ws.Cells.Item(24, 17).Value = 490.6740854295783  # Q24
$ws.Cells.Item(24, 18).Value = 2944.04451257747  # R24
$ws.Cells.Item(24, 19).Value = 0.010697679266424433  # S24
$ws.Cells.Item(24, 20).Value = 0.010964002260210503  # T24
$ws.Cells.Item(25, 7).Value = 23.049115  # G25
$ws.Cells.Item(25, 8).Value = 46.09823  # H25
$ws.Cells.Item(25, 9).Value = 0.2006333249297675  # I25
$ws.Cells.Item(25, 10).Value = 0.16498447796139157  # J25
$ws.Cells.Item(25, 13).Value = 102.442257  # M25
$ws.Cells.Item(25, 14).Value = 307.326771  # N25
$ws.Cells.Item(25, 15).Value = 0.25658234727063134  # O25
$ws.Cells.Item(25, 16).Value = 0.31979104694013644  # P25
$ws.Cells.Item(25, 17).Value = 2361.203362452555  # Q25
$ws.Cells.Item(25, 18).Value = 14167.22017471533  # R25
$ws.Cells.Item(25, 19).Value = 0.051478969451191015  # S25
$ws.Cells.Item(25, 20).Value = 0.05276055893614528  # T25
$ws.Cells.Item(26, 7).Value = 0.212799  # G26
$ws.Cells.Item(26, 8).Value = 0.6383969999999999  # H26
$ws.Cells.Item(26, 9).Value = 0.0018523301615584627  # I26
$ws.Cells.Item(26, 10).Value = 0.002284807806658053  # J26
$ws.Cells.Item(26, 13).Value = 228.691925  # M26
$ws.Cells.Item(26, 14).Value = 457.38385  # N26
$ws.Cells.Item(26, 15).Value = 0.5727940074410814  # O26
$ws.Cells.Item(26, 16).Value = 0.475934002654817  # P26
$ws.Cells.Item(26, 17).Value = 48.665412948074994  # Q26
$ws.Cells.Item(26, 18).Value = 291.9924776884499  # R26
$ws.Cells.Item(26, 19).Value = 0.0010610036163430576  # S26
$ws.Cells.Item(26, 20).Value = 0.0010874177247197404  # T26
$ws.Cells.Item(27, 7).Value = 0.212799  # G27
$ws.Cells.Item(27, 8).Value = 0.6383969999999999  # H27
$ws.Cells.Item(27, 9).Value = 0.0018523301615584627  # I27
$ws.Cells.Item(27, 10).Value = 0.002284807806658053  # J27
$ws.Cells.Item(27, 15).Value = 0.0011895817926618896  # O27
$ws.Cells.Item(27, 16).Value = 0.0014826335909033635  # P27
$ws.Cells.Item(27, 17).Value = 0.101068601318  # Q27
$ws.Cells.Item(27, 18).Value = 0.9096174118619998  # R27
$ws.Cells.Item(27, 19).Value = 0.0000022034982341884037  # S27
$ws.Cells.Item(27, 20).Value = 0.0000033875328029094673  # T27
$ws.Cells.Item(28, 7).Value = 0.212799  # G28
$ws.Cells.Item(28, 8).Value = 0.6383969999999999  # H28
$ws.Cells.Item(28, 9).Value = 0.0018523301615584627  # I28
$ws.Cells.Item(28, 10).Value = 0.002284807806658053  # J28
$ws.Cells.Item(28, 13).Value = 38.30460633333333  # M28
$ws.Cells.Item(28, 14).Value = 114.913819  # N28
$ws.Cells.Item(28, 15).Value = 0.0959397624779407  # O28
$ws.Cells.Item(28, 16).Value = 0.11957438776428411  # P28
$ws.Cells.Item(28, 17).Value = 8.151181923127  # Q28
$ws.Cells.Item(28, 18).Value = 73.36063730814298  # R28
$ws.Cells.Item(28, 19).Value = 0.00017771211573064442  # S28
$ws.Cells.Item(28, 20).Value = 0.00027320449464019353  # T28
$ws.Cells.Item(29, 7).Value = 0.212799  # G29
$ws.Cells.Item(29, 8).Value = 0.6383969999999999  # H29
$ws.Cells.Item(29, 9).Value = 0.0018523301615584627  # I29
$ws.Cells.Item(29, 10).Value = 0.002284807806658053  # J29
$ws.Cells.Item(29, 13).Value = 8.054906  # M29
$ws.Cells.Item(29, 14).Value = 16.109812  # N29
$ws.Cells.Item(29, 15).Value = 0.02017474769737153  # O29
$ws.Cells.Item(29, 16).Value = 0.01676317890799293  # P29
$ws.Cells.Item(29, 17).Value = 1.7140759418940001  # Q29
$ws.Cells.Item(29, 18).Value = 10.284455651363999  # R29
$ws.Cells.Item(29, 19).Value = 0.000037370293661673425  # S29
$ws.Cells.Item(29, 20).Value = 0.00003830064203338787  # T29
$ws.Cells.Item(30, 7).Value = 0.212799  # G30
$ws.Cells.Item(30, 8).Value = 0.6383969999999999  # H30
$ws.Cells.Item(30, 9).Value = 0.0018523301615584627  # I30
$ws.Cells.Item(30, 10).Value = 0.002284807806658053  # J30
$ws.Cells.Item(30, 13).Value = 21.28819633333333  # M30
$ws.Cells.Item(30, 14).Value = 63.864589  # N30
$ws.Cells.Item(30, 15).Value = 0.05331955332031306  # O30
$ws.Cells.Item(30, 16).Value = 0.06645475014186618  # P30
$ws.Cells.Item(30, 17).Value = 4.530106891537  # Q30
$ws.Cells.Item(30, 18).Value = 40.77096202383299  # R30
$ws.Cells.Item(30, 19).Value = 0.00009876541681604055  # S30
$ws.Cells.Item(30, 20).Value = 0.0001518363319136462  # T30
$ws.Cells.Item(31, 7).Value = 0.212799  # G31
$ws.Cells.Item(31, 8).Value = 0.6383969999999999  # H31
$ws.Cells.Item(31, 9).Value = 0.0018523301615584627  # I31
$ws.Cells.Item(31, 10).Value = 0.002284807806658053  # J31
$ws.Cells.Item(31, 13).Value = 102.442257  # M31
$ws.Cells.Item(31, 14).Value = 307.326771  # N31
$ws.Cells.Item(31, 15).Value = 0.25658234727063134  # O31
$ws.Cells.Item(31, 16).Value = 0.31979104694013644  # P31
$ws.Cells.Item(31, 17).Value = 21.799609847343  # Q31
$ws.Cells.Item(31, 18).Value = 196.19648862608696  # R31
$ws.Cells.Item(31, 19).Value = 0.0004752752207728581  # S31
$ws.Cells.Item(31, 20).Value = 0.0007306610805481756  # T31
$ws.Cells.Item(32, 5).Value = 1.0  # E32
$ws.Cells.Item(32, 6).Value = 0.3333333333333333  # F32
$ws.Cells.Item(32, 7).Value = 0.05774266666666666  # G32
$ws.Cells.Item(32, 8).Value = 0.173228  # H32
$ws.Cells.Item(32, 9).Value = 0.0005026268125107878  # I32
$ws.Cells.Item(32, 10).Value = 0.0006199789264858095  # J32
$ws.Cells.Item(32, 13).Value = 228.691925  # M32
$ws.Cells.Item(32, 14).Value = 457.38385  # N32
$ws.Cells.Item(32, 15).Value = 0.5727940074410814  # O32
$ws.Cells.Item(32, 16).Value = 0.475934002654817  # P32
$ws.Cells.Item(32, 17).Value = 13.20528159463333  # Q32
$ws.Cells.Item(32, 18).Value = 79.23168956779999  # R32
$ws.Cells.Item(32, 19).Value = 0.0002879016261853912  # S32
$ws.Cells.Item(32, 20).Value = 0.00029506905204402785  # T32
$ws.Cells.Item(33, 5).Value = 1.0  # E33
$ws.Cells.Item(33, 6).Value = 0.3333333333333333  # F33
$ws.Cells.Item(33, 7).Value = 0.05774266666666666  # G33
$ws.Cells.Item(33, 8).Value = 0.173228  # H33
$ws.Cells.Item(33, 9).Value = 0.0005026268125107878  # I33
$ws.Cells.Item(33, 10).Value = 0.0006199789264858095  # J33
$ws.Cells.Item(33, 15).Value = 0.0011895817926618896  # O33
$ws.Cells.Item(33, 16).Value = 0.0014826335909033635  # P33
$ws.Cells.Item(33, 17).Value = 0.027424802543111106  # Q33
$ws.Cells.Item(33, 18).Value = 0.246823222888  # R33
$ws.Cells.Item(33, 19).Value = 0.0000005979157046665144  # S33
$ws.Cells.Item(33, 20).Value = 0.0000009192015820600681  # T33
$ws.Cells.Item(34, 5).Value = 1.0  # E34
$ws.Cells.Item(34, 6).Value = 0.3333333333333333  # F34
$ws.Cells.Item(34, 7).Value = 0.05774266666666666  # G34
$ws.Cells.Item(34, 8).Value = 0.173228  # H34
$ws.Cells.Item(34, 9).Value = 0.0005026268125107878  # I34
$ws.Cells.Item(34, 10).Value = 0.0006199789264858095  # J34
$ws.Cells.Item(34, 13).Value = 38.30460633333333  # M34
$ws.Cells.Item(34, 14).Value = 114.913819  # N34
$ws.Cells.Item(34, 15).Value = 0.0959397624779407  # O34
$ws.Cells.Item(34, 16).Value = 0.11957438776428411  # P34
$ws.Cells.Item(34, 17).Value = 2.211810115303555  # Q34
$ws.Cells.Item(34, 18).Value = 19.906291037732  # R34
$ws.Cells.Item(34, 19).Value = 0.00004822189700732941  # S34
$ws.Cells.Item(34, 20).Value = 0.00007413360056129877  # T34
$ws.Cells.Item(35, 5).Value = 1.0  # E35
$ws.Cells.Item(35, 6).Value = 0.3333333333333333  # F35
$ws.Cells.Item(35, 7).Value = 0.05774266666666666  # G35
$ws.Cells.Item(35, 8).Value = 0.173228  # H35
$ws.Cells.Item(35, 9).Value = 0.0005026268125107878  # I35
$ws.Cells.Item(35, 10).Value = 0.0006199789264858095  # J35
$ws.Cells.Item(35, 13).Value = 8.054906  # M35
$ws.Cells.Item(35, 14).Value = 16.109812  # N35
$ws.Cells.Item(35, 15).Value = 0.02017474769737153  # O35
$ws.Cells.Item(35, 16).Value = 0.01676317890799293  # P35
$ws.Cells.Item(35, 17).Value = 0.4651117521893333  # Q35
$ws.Cells.Item(35, 18).Value = 2.7906705131360003  # R35
$ws.Cells.Item(35, 19).Value = 0.000010140369128339206  # S35
$ws.Cells.Item(35, 20).Value = 0.000010392817663867022  # T35
$ws.Cells.Item(36, 5).Value = 1.0  # E36
$ws.Cells.Item(36, 6).Value = 0.3333333333333333  # F36
$ws.Cells.Item(36, 7).Value = 0.05774266666666666  # G36
$ws.Cells.Item(36, 8).Value = 0.173228  # H36
$ws.Cells.Item(36, 9).Value = 0.0005026268125107878  # I36
$ws.Cells.Item(36, 10).Value = 0.0006199789264858095  # J36
$ws.Cells.Item(36, 13).Value = 21.28819633333333  # M36
$ws.Cells.Item(36, 14).Value = 63.864589  # N36
$ws.Cells.Item(36, 15).Value = 0.05331955332031306  # O36
$ws.Cells.Item(36, 16).Value = 0.06645475014186618  # P36
$ws.Cells.Item(36, 17).Value = 1.229237224810222  # Q36
$ws.Cells.Item(36, 18).Value = 11.063135023292  # R36
$ws.Cells.Item(36, 19).Value = 0.000026799837129887944  # S36
$ws.Cells.Item(36, 20).Value = 0.00004120054465283689  # T36
$ws.Cells.Item(37, 5).Value = 1.0  # E37
$ws.Cells.Item(37, 6).Value = 0.3333333333333333  # F37
$ws.Cells.Item(37, 7).Value = 0.05774266666666666  # G37
$ws.Cells.Item(37, 8).Value = 0.173228  # H37
$ws.Cells.Item(37, 9).Value = 0.0005026268125107878  # I37
$ws.Cells.Item(37, 10).Value = 0.0006199789264858095  # J37
$ws.Cells.Item(37, 13).Value = 102.442257  # M37
$ws.Cells.Item(37, 14).Value = 307.326771  # N37
$ws.Cells.Item(37, 15).Value = 0.25658234727063134  # O37
$ws.Cells.Item(37, 16).Value = 0.31979104694013644  # P37
$ws.Cells.Item(37, 17).Value = 5.915289098531999  # Q37
$ws.Cells.Item(37, 18).Value = 53.237601886788  # R37
$ws.Cells.Item(37, 19).Value = 0.00012896516735517346  # S37
$ws.Cells.Item(37, 20).Value = 0.0001982637099817189  # T37
